$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.526.46"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "2.989.41"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'381.57"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'104.22"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").Value = "'36.73"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").Value = "'0.137"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "3.465.56"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "'18.51"
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "'7.80"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "2.983.33"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "'11.25"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "'0.994"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "51.604.08"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "'12.51"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").Value = "'70.40"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "'267.51"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "'3.22"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "'8.09"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +4.79%  "
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'26.10"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'10.41"
$ws.Range("E32").Value = "  +4.25%  "
$ws.Range("D33").Value = "'34.65"
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("D34").Value = "'51.39"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +6.15%  "
$ws.Range("D39").Value = "'17.07"
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("D40").Value = "'2.59"
$ws.Range("E40").Value = "  +5.65%  "
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.86"
$ws.Range("E43").Value = "  +15.00%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'122.75"
$ws.Range("E44").Value = "  +4.63%  "
$ws.Range("D45").Value = "'21.43"
$ws.Range("E45").Value = "  +1.89%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "'0.271"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "2.039.81"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").Value = "3.286.44"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("E51").Value = "  +2.92%  "
